$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Reorder worksheet tabs.
#    Before: Companies, AppName, ModuleName, Users, CoverageOfficer, TabName
#    After : Users, AppName, ModuleName, Companies, CoverageOfficer, TabName
# ---------------------------------------------------------------------------
$wsCompanies = $wb.Worksheets.Item("Companies")
$wsUsers     = $wb.Worksheets.Item("Users")
$wsCompanies.Move($wsUsers)
# Now order is: AppName, ModuleName, Companies, Users, CoverageOfficer, TabName
$wsAppName = $wb.Worksheets.Item("AppName")
$wsUsers   = $wb.Worksheets.Item("Users")
$wsUsers.Move($wsAppName)
# Now order is: Users, AppName, ModuleName, Companies, CoverageOfficer, TabName

$wsCompanies = $wb.Worksheets.Item("Companies")
$wsUsers     = $wb.Worksheets.Item("Users")
$wsCoverage  = $wb.Worksheets.Item("CoverageOfficer")
$wsTabName   = $wb.Worksheets.Item("TabName")

# ---------------------------------------------------------------------------
# 2. Stash the "text / left-top aligned" cell format (used on a couple of
#    blank helper cells) onto a scratch cell before the sheets carrying it
#    get cleared, so it can be re-applied later without creating a brand
#    new style entry.
# ---------------------------------------------------------------------------
$wsCoverage.Range("C3").Copy()
$wsTabName.Range("Z100").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Clear out the old content of the three sheets being rewritten.
# ---------------------------------------------------------------------------
$wsCompanies.UsedRange.Clear()
$wsUsers.UsedRange.Clear()
$wsCoverage.UsedRange.Clear()

# ---------------------------------------------------------------------------
# 4. Write the new cell values. The order below matters: it controls the
#    order in which brand-new shared strings are appended.
# ---------------------------------------------------------------------------

# -- Companies --------------------------------------------------------------
$wsCompanies.Range("A1:H1").Font.Bold = $true
$wsCompanies.Range("A1").Value = "CompaniesName"
$wsCompanies.Range("B1").Value = "CompanyType"

$wsCompanies.Range("A2").Value = "ADK Holdings"
$wsCompanies.Range("B2").Value = "Operating Company"
$wsCompanies.Range("D2").Value = "Varta AG"

$wsCompanies.Range("A3").Value = "H.I.G. Capital Management, Inc."
$wsCompanies.Range("B3").Value = "Capital Provider"
$wsCompanies.Range("F3").Value = "Audax Management Company, LLC"

# -- CoverageOfficer (row 3 first, then row 2, to match shared-string order) -
$wsCoverage.Range("A1:H1").Font.Bold = $true
$wsCoverage.Range("A1").Value = "CoverageOfficer"

$wsCoverage.Range("A3").Value = "Daniel Gossels"
$wsCoverage.Range("B3").Value = "Jim Lavelle"
$wsCoverage.Range("D3").Value = "Michael Morabito"

# -- Users --------------------------------------------------------------
$wsUsers.Range("A1:B1").Font.Bold = $true
$wsUsers.Range("A1").Value = "StdUser"

$wsUsers.Range("A2").Value = "James Craven"

$wsUsers.Range("A3").Value = "Ajay Nair"
$wsUsers.Range("B3").Value = "System Admin"

# -- CoverageOfficer row 2 (after Users, to match shared-string order) ------
$wsCoverage.Range("A2").Value = "Tomohiko Kubota"
$wsCoverage.Range("B2").Value = "Steve Hughes"
$wsCoverage.Range("C2").Value = "Yuta Nakamura"

# ---------------------------------------------------------------------------
# 5. Re-apply the text/left-top alignment style to the two blank formatted
#    cells on the "Companies" sheet, then drop the scratch cell.
# ---------------------------------------------------------------------------
$wsTabName.Range("Z100").Copy()
$wsCompanies.Range("A3").PasteSpecial(-4122)
$wsCompanies.Range("G3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsTabName.Range("Z100").Clear()

# ---------------------------------------------------------------------------
# 6. Selections (set non-active sheets first, active sheet/tab last).
# ---------------------------------------------------------------------------
$wsUsers.Range("B37").Select()
$wsCompanies.Range("D30").Select()
$wsTabName.Range("K24").Select()

$wsCoverage.Activate()
$wsCoverage.Range("H20").Select()
